$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("YDS")

$ws.Range("B2").Value = "NCT(2.866960236765078, 1.2108580379941865, 0.09299699794564023, 2.5375334989226905)"
$ws.Range("C2").Value = "NIG(1.2874689748737906, 0.9180171389018468, 4.050287401345738, 6.238441027794778)"
$ws.Range("D2").Value = "NIG(0.8203995437981877, 0.5203011854978306, 1.7059292795422563, 3.07632641683195)"
$ws.Range("E2").Value = "JSU(-1.4571164103794467, 1.3376898055101358, 2.057926571428432, 5.448513960767446)"
